# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 10:46"

# A Coruña (row 16) - Muertes 5 -> 12
$ws.Range("E16").Value = 12

# Pontevedra (row 18) - Muertes 2 -> 3
$ws.Range("E18").Value = 3

# Mallorca (row 31) - Recuperados 0 -> 197, Muertes 8 -> 9
$ws.Range("D31").Value = 197
$ws.Range("E31").Value = 9

# Ourense (row 44) - Muertes 0 -> 2
$ws.Range("E44").Value = 2

# Lugo (row 47) - Muertes 1 -> 2
$ws.Range("E47").Value = 2

# Ibiza (row 54) - Recuperados 0 -> 20, Muertes 8 -> 1
$ws.Range("D54").Value = 20
$ws.Range("E54").Value = 1

# Menorca (row 55) - Recuperados 0 -> 13, Muertes 8 -> 0
$ws.Range("D55").Value = 13
$ws.Range("E55").Value = 0
